$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the contents of row 2 and row 3 for the columns that differ
# (A, B, E, F, G, I, M). Columns C, D, H, K, L, N etc. are identical
# between the two rows, so no visible change occurs there.

$cols = @("A", "B", "E", "F", "G", "I", "M")

foreach ($col in $cols) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")

    $v2 = $cell2.Value()
    $v3 = $cell3.Value()

    $cell2.Value = $v3
    $cell3.Value = $v2
}
